$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2869.2307
$ws.Range("J64").Value = 2757.1428
$ws.Range("L64").Value = 2757.1428
$ws.Range("N64").Value = -3253.1428
$ws.Range("H67").Value = 2869.2307
$ws.Range("J67").Value = 2757.1428
$ws.Range("L67").Value = 2757.1428
$ws.Range("N67").Value = -4473.1428
$ws.Range("H76").Value = 3315.0303
$ws.Range("I76").Value = 3269.3076
$ws.Range("J76").Value = 3484.8572
$ws.Range("K76").Value = 3269.3076
$ws.Range("L76").Value = 3484.8572
$ws.Range("M76").Value = -2954.3076
$ws.Range("N76").Value = -4114.8572
$ws.Range("H79").Value = 3315.0303
$ws.Range("I79").Value = 3269.3076
$ws.Range("J79").Value = 3484.8572
$ws.Range("K79").Value = 3269.3076
$ws.Range("L79").Value = 3484.8572
$ws.Range("M79").Value = -2177.3076
$ws.Range("N79").Value = -5668.8572
$ws.Range("H129").Value = 925.6389
$ws.Range("J129").Value = 959.5146999999999
$ws.Range("L129").Value = 2878.5441
$ws.Range("N129").Value = -12878.5441
$ws.Range("H138").Value = 3064.4795
$ws.Range("I138").Value = 1828.421
$ws.Range("J138").Value = 3361.7595
$ws.Range("K138").Value = 5485.263
$ws.Range("L138").Value = 10085.2785
$ws.Range("M138").Value = -345.2629999999999
$ws.Range("N138").Value = -20365.2785
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6494.7397
$ws.Range("I32").Value = 3652.1738
$ws.Range("J32").Value = 13759.074
$ws.Range("K32").Value = 3652.1738
$ws.Range("L32").Value = 13759.074
$ws.Range("M32").Value = -3365.1738
$ws.Range("N32").Value = -14333.074
$ws.Range("H45").Value = 1620.2
$ws.Range("I45").Value = 1002.8
$ws.Range("J45").Value = 2237.6
$ws.Range("K45").Value = 1002.8
$ws.Range("L45").Value = 2237.6
$ws.Range("M45").Value = -625.8
$ws.Range("N45").Value = -2991.6
$ws.Range("H63").Value = 13853421
$ws.Range("J63").Value = 3500
$ws.Range("L63").Value = 3500
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 13853421
$ws.Range("J66").Value = 3500
$ws.Range("L66").Value = 17500
$ws.Range("N66").Value = -24364
$ws.Range("H74").Value = 1325.2449
$ws.Range("I74").Value = 886.64105
$ws.Range("J74").Value = 3035.8
$ws.Range("K74").Value = 886.64105
$ws.Range("L74").Value = 3035.8
$ws.Range("M74").Value = -12.64104999999995
$ws.Range("N74").Value = -4783.8
$ws.Range("H77").Value = 1325.2449
$ws.Range("I77").Value = 886.64105
$ws.Range("J77").Value = 3035.8
$ws.Range("K77").Value = 4433.20525
$ws.Range("L77").Value = 15179
$ws.Range("M77").Value = -65.20524999999998
$ws.Range("N77").Value = -23915
$ws.Range("H110").Value = 609.8
$ws.Range("I110").Value = 566.44446
$ws.Range("K110").Value = 566.44446
$ws.Range("M110").Value = 1478.55554
$ws.Range("H134").Value = 48412.715
$ws.Range("J134").Value = 48412.715
$ws.Range("L134").Value = 48412.715
$ws.Range("N134").Value = -58552.715
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 799.8333
$ws.Range("I5").Value = 299.5
$ws.Range("J5").Value = 1050
$ws.Range("K5").Value = 299.5
$ws.Range("L5").Value = 1050
$ws.Range("M5").Value = -187.5
$ws.Range("N5").Value = -1274
$ws.Range("H31").Value = 2753.2083
$ws.Range("I31").Value = 1053.1515
$ws.Range("J31").Value = 6493.3335
$ws.Range("K31").Value = 1053.1515
$ws.Range("L31").Value = 6493.3335
$ws.Range("M31").Value = -758.1514999999999
$ws.Range("N31").Value = -7083.3335
$ws.Range("H34").Value = 2753.2083
$ws.Range("I34").Value = 1053.1515
$ws.Range("J34").Value = 6493.3335
$ws.Range("K34").Value = 1053.1515
$ws.Range("L34").Value = 6493.3335
$ws.Range("M34").Value = -851.1514999999999
$ws.Range("N34").Value = -6897.3335
$ws.Range("H58").Value = 1821.7428
$ws.Range("I58").Value = 1558.5374
$ws.Range("J58").Value = 7700
$ws.Range("K58").Value = 1558.5374
$ws.Range("L58").Value = 7700
$ws.Range("M58").Value = -1355.5374
$ws.Range("N58").Value = -8106
$ws.Range("H132").Value = 2105.3774
$ws.Range("I132").Value = 1666.8223
$ws.Range("J132").Value = 4572.25
$ws.Range("K132").Value = 5000.4669
$ws.Range("L132").Value = 13716.75
$ws.Range("M132").Value = -2470.4669
$ws.Range("N132").Value = -18776.75
$ws.Range("H134").Value = 3100.3845
$ws.Range("I134").Value = 3119.9607
$ws.Range("J134").Value = 3029.0715
$ws.Range("K134").Value = 9359.882100000001
$ws.Range("L134").Value = 9087.2145
$ws.Range("M134").Value = -6824.882100000001
$ws.Range("N134").Value = -14157.2145
$ws.Range("H136").Value = 1821.7428
$ws.Range("I136").Value = 1558.5374
$ws.Range("J136").Value = 7700
$ws.Range("K136").Value = 4675.6122
$ws.Range("L136").Value = 23100
$ws.Range("M136").Value = -2125.6122
$ws.Range("N136").Value = -28200
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2381031
$ws.Range("I2").Value = 56
$ws.Range("K2").Value = 336
$ws.Range("M2").Value = -223
$ws.Range("H9").Value = 334000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 334000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1002000
$ws.Range("N9").Value = -1002448
$ws.Range("M9").ClearContents()
$ws.Range("H10").Value = 3200
$ws.Range("J10").Value = 3200
$ws.Range("L10").Value = 9600
$ws.Range("N10").Value = -9878
$ws.Range("H133").Value = 2921.25
$ws.Range("I133").Value = 2565
$ws.Range("J133").Value = 3990
$ws.Range("K133").Value = 7695
$ws.Range("L133").Value = 11970
$ws.Range("M133").Value = -2635
$ws.Range("N133").Value = -22090
$ws.Range("H137").Value = 6458.6665
$ws.Range("I137").Value = 2622.1428
$ws.Range("J137").Value = 10590.308
$ws.Range("K137").Value = 7866.428400000001
$ws.Range("L137").Value = 31770.924
$ws.Range("M137").Value = -2766.428400000001
$ws.Range("N137").Value = -41970.924
$ws.Range("H140").Value = 74057.14
$ws.Range("I140").Value = 168800
$ws.Range("K140").Value = 506400
$ws.Range("M140").Value = -501220
$ws.Range("H141").Value = 8329.923000000001
$ws.Range("J141").Value = 8320
$ws.Range("L141").Value = 24960
$ws.Range("N141").Value = -35320
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5589.9443
$ws.Range("I70").Value = 5316.95
$ws.Range("K70").Value = 5316.95
$ws.Range("M70").Value = -5046.95
$ws.Range("H73").Value = 5589.9443
$ws.Range("I73").Value = 5316.95
$ws.Range("K73").Value = 5316.95
$ws.Range("M73").Value = -4380.95
$ws.Range("H80").Value = 31252574
$ws.Range("I80").Value = 250000000
$ws.Range("J80").Value = 2941.4285
$ws.Range("K80").Value = 250000000
$ws.Range("L80").Value = 2941.4285
$ws.Range("M80").Value = -249999002
$ws.Range("N80").Value = -4937.4285
$ws.Range("H83").Value = 31252574
$ws.Range("I83").Value = 250000000
$ws.Range("J83").Value = 2941.4285
$ws.Range("K83").Value = 1250000000
$ws.Range("L83").Value = 14707.1425
$ws.Range("M83").Value = -1249995008
$ws.Range("N83").Value = -24691.1425
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2041.2963
$ws.Range("I46").Value = 1689.1111
$ws.Range("J46").Value = 2217.389
$ws.Range("K46").Value = 1689.1111
$ws.Range("L46").Value = 2217.389
$ws.Range("M46").Value = -1501.1111
$ws.Range("N46").Value = -2593.389
$ws.Range("H69").Value = 359600
$ws.Range("J69").Value = 359600
$ws.Range("L69").Value = 359600
$ws.Range("N69").Value = -361222
$ws.Range("H72").Value = 359600
$ws.Range("J72").Value = 359600
$ws.Range("L72").Value = 1078800
$ws.Range("N72").Value = -1086912
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6948108.5
$ws.Range("I132").Value = 4054.7585
$ws.Range("K132").Value = 12164.2755
$ws.Range("M132").Value = -9634.2755
